$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.867.73"
$ws.Cells.Item(2, 5).Value = "  -0.28%  "
$ws.Cells.Item(3, 4).Value = "1.861.42"
$ws.Cells.Item(3, 5).Value = "  +0.05%  "
$ws.Cells.Item(4, 4).Value = "'0.9992"
$ws.Cells.Item(5, 4).Value = "'304.72"
$ws.Cells.Item(5, 5).Value = "  -0.38%  "
$ws.Cells.Item(6, 4).Value = "'0.9996"
$ws.Cells.Item(6, 5).Value = "  -0.13%  "
$ws.Cells.Item(7, 4).Value = "'0.5047"
$ws.Cells.Item(7, 5).Value = "  -0.31%  "
$ws.Cells.Item(8, 4).Value = "'0.3643"
$ws.Cells.Item(8, 5).Value = "  -2.33%  "
$ws.Cells.Item(9, 5).Value = "  +0.54%  "
$ws.Cells.Item(10, 4).Value = "'0.8916"
$ws.Cells.Item(10, 5).Value = "  +0.53%  "
$ws.Cells.Item(11, 4).Value = "'20.66"
$ws.Cells.Item(11, 5).Value = "  +0.70%  "
$ws.Cells.Item(12, 2).Value = "WrappedEther"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(12, 4).Value = "1.866.18"
$ws.Cells.Item(12, 5).Value = "  +0.29%  "
$ws.Cells.Item(13, 2).Value = "TRON"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(13, 4).Value = "'0.07504"
$ws.Cells.Item(13, 5).Value = "  -0.57%  "
$ws.Cells.Item(14, 4).Value = "'94.91"
$ws.Cells.Item(14, 5).Value = "  +6.78%  "
$ws.Cells.Item(15, 4).Value = "'5.226"
$ws.Cells.Item(15, 5).Value = "  -1.16%  "
$ws.Cells.Item(16, 4).Value = "'0.9997"
$ws.Cells.Item(16, 5).Value = "  -0.13%  "
$ws.Cells.Item(17, 4).Value = "'0.000008515"
$ws.Cells.Item(17, 5).Value = "  +1.83%  "
$ws.Cells.Item(18, 5).Value = "  +1.13%  "
$ws.Cells.Item(19, 4).Value = "'0.9996"
$ws.Cells.Item(19, 5).Value = "  -0.13%  "
$ws.Cells.Item(20, 4).Value = "26.923.23"
$ws.Cells.Item(20, 5).Value = "  -0.30%  "
$ws.Cells.Item(21, 4).Value = "'5.025"
$ws.Cells.Item(21, 5).Value = "  -0.35%  "
$ws.Cells.Item(22, 4).Value = "2.103.85"
$ws.Cells.Item(22, 5).Value = "  +0.72%  "
$ws.Cells.Item(23, 5).Value = "  -1.00%  "
$ws.Cells.Item(24, 5).Value = "  -0.89%  "
$ws.Cells.Item(25, 4).Value = "'147.76"
$ws.Cells.Item(25, 5).Value = "  +0.32%  "
$ws.Cells.Item(26, 5).Value = "  -3.37%  "
$ws.Cells.Item(27, 5).Value = "  -0.47%  "
$ws.Cells.Item(28, 4).Value = "'2.083"
$ws.Cells.Item(28, 5).Value = "  -0.21%  "
$ws.Cells.Item(29, 4).Value = "'113.23"
$ws.Cells.Item(30, 4).Value = "'4.705"
$ws.Cells.Item(30, 5).Value = "  +1.21%  "
$ws.Cells.Item(31, 4).Value = "'4.667"
$ws.Cells.Item(31, 5).Value = "  +0.53%  "
$ws.Cells.Item(32, 5).Value = "  +1.93%  "
$ws.Cells.Item(33, 5).Value = "  +0.62%  "
$ws.Cells.Item(34, 4).Value = "'0.7502"
$ws.Cells.Item(34, 5).Value = "  +3.17%  "
$ws.Cells.Item(35, 4).Value = "'2.958"
$ws.Cells.Item(35, 5).Value = "  -3.02%  "
$ws.Cells.Item(36, 4).Value = "'1.152"
$ws.Cells.Item(36, 5).Value = "  +0.36%  "
$ws.Cells.Item(37, 4).Value = "'3.250"
$ws.Cells.Item(37, 5).Value = "  +7.07%  "
$ws.Cells.Item(38, 4).Value = "'2.571"
$ws.Cells.Item(38, 5).Value = "  +5.41%  "
$ws.Cells.Item(39, 4).Value = "'0.02001"
$ws.Cells.Item(39, 5).Value = "  -1.80%  "
$ws.Cells.Item(40, 4).Value = "'0.5575"
$ws.Cells.Item(40, 5).Value = "  +5.00%  "
$ws.Cells.Item(41, 5).Value = "  -0.07%  "
$ws.Cells.Item(42, 4).Value = "'6.557"
$ws.Cells.Item(42, 5).Value = "  -0.17%  "
$ws.Cells.Item(43, 4).Value = "'116.06"
$ws.Cells.Item(43, 5).Value = "  +0.68%  "
$ws.Cells.Item(44, 4).Value = "'8.564"
$ws.Cells.Item(44, 5).Value = "  +3.65%  "
$ws.Cells.Item(45, 4).Value = "'0.1470"
$ws.Cells.Item(45, 5).Value = "  +0.10%  "
$ws.Cells.Item(46, 5).Value = "  +2.07%  "
$ws.Cells.Item(47, 4).Value = "'0.9993"
$ws.Cells.Item(47, 5).Value = "  -0.16%  "
$ws.Cells.Item(48, 4).Value = "'10.07"
$ws.Cells.Item(48, 5).Value = "  +1.45%  "
$ws.Cells.Item(49, 4).Value = "'1.562"
$ws.Cells.Item(49, 5).Value = "  +0.51%  "
$ws.Cells.Item(50, 4).Value = "'36.70"
$ws.Cells.Item(50, 5).Value = "  +0.65%  "
$ws.Cells.Item(51, 4).Value = "'62.92"
$ws.Cells.Item(51, 5).Value = "  -1.64%  "
